$d = $word.ActiveDocument

# Update the date/weekday heading.
$found = $d.Content.Find.Execute("2024-10-09 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-10-10 Thursday", 2)

# Update the division-problem answer table (first table in the document).
# The table has data rows at 1, 5, 9, 13, 17 (1-based), each with 5 columns.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "585÷4=146, 1"
$t.Cell(1, 2).Range.Text  = "485÷4=121, 1"
$t.Cell(1, 3).Range.Text  = "587÷6=97, 5"
$t.Cell(1, 4).Range.Text  = "818÷6=136, 2"
$t.Cell(1, 5).Range.Text  = "547÷8=68, 3"

$t.Cell(5, 1).Range.Text  = "766÷8=95, 6"
$t.Cell(5, 2).Range.Text  = "537÷8=67, 1"
$t.Cell(5, 3).Range.Text  = "933÷4=233, 1"
$t.Cell(5, 4).Range.Text  = "211÷6=35, 1"
$t.Cell(5, 5).Range.Text  = "820÷5=164, 0"

$t.Cell(9, 1).Range.Text  = "268÷8=33, 4"
$t.Cell(9, 2).Range.Text  = "281÷9=31, 2"
$t.Cell(9, 3).Range.Text  = "383÷5=76, 3"
$t.Cell(9, 4).Range.Text  = "475÷6=79, 1"
$t.Cell(9, 5).Range.Text  = "494÷8=61, 6"

$t.Cell(13, 1).Range.Text = "996÷6=166, 0"
$t.Cell(13, 2).Range.Text = "653÷6=108, 5"
$t.Cell(13, 3).Range.Text = "331÷2=165, 1"
$t.Cell(13, 4).Range.Text = "860÷7=122, 6"
$t.Cell(13, 5).Range.Text = "486÷7=69, 3"

$t.Cell(17, 1).Range.Text = "677÷4=169, 1"
$t.Cell(17, 2).Range.Text = "132÷6=22, 0"
$t.Cell(17, 3).Range.Text = "300÷4=75, 0"
$t.Cell(17, 4).Range.Text = "961÷9=106, 7"
$t.Cell(17, 5).Range.Text = "521÷9=57, 8"

Write-Output "done"
